$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.794.08'
$ws.Range("E2").Value = '  -4.04%  '
$ws.Range("D3").Value = '3.386.38'
$ws.Range("E3").Value = '  -4.52%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'563.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.17%  '
$ws.Range("D6").Value = "'184.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.99%  '
$ws.Range("E7").Value = '  -2.01%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '3.377.55'
$ws.Range("E9").Value = '  -4.55%  '
$ws.Range("E10").Value = '  -8.47%  '
$ws.Range("E11").Value = '  -4.84%  '
$ws.Range("D12").Value = "'48.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.28%  '
$ws.Range("E13").Value = '  -6.58%  '
$ws.Range("E14").Value = '  -5.81%  '
$ws.Range("D15").Value = '3.921.88'
$ws.Range("E15").Value = '  -4.55%  '
$ws.Range("D16").Value = "'610.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -11.55%  '
$ws.Range("D17").Value = "'18.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.92%  '
$ws.Range("D18").Value = '66.685.42'
$ws.Range("E18").Value = '  -4.36%  '
$ws.Range("D19").Value = '3.386.56'
$ws.Range("E19").Value = '  -4.81%  '
$ws.Range("E20").Value = '  -2.91%  '
$ws.Range("D21").Value = "'11.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.43%  '
$ws.Range("D22").Value = "'0.923"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.01%  '
$ws.Range("D23").Value = "'17.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.12%  '
$ws.Range("D24").Value = "'5.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.62%  '
$ws.Range("D25").Value = "'99.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.45%  '
$ws.Range("D26").Value = "'4.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.02%  '
$ws.Range("D27").Value = "'6.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.69%  '
$ws.Range("E28").Value = '  -7.39%  '
$ws.Range("D29").Value = "'9.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.89%  '
$ws.Range("E30").Value = '  -8.56%  '
$ws.Range("D31").Value = "'31.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.06%  '
$ws.Range("D32").Value = "'3.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -10.86%  '
$ws.Range("D33").Value = "'6.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.77%  '
$ws.Range("D34").Value = "'11.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.91%  '
$ws.Range("D35").Value = "'559.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.70%  '
$ws.Range("D36").Value = '3.904.05'
$ws.Range("E36").Value = '  +2.46%  '
$ws.Range("E37").Value = '  -4.82%  '
$ws.Range("D38").Value = "'58.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.24%  '
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.18%  '
$ws.Range("D40").Value = "'3.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.67%  '
$ws.Range("E41").Value = '  +25.91%  '
$ws.Range("E42").Value = '  -11.35%  '
$ws.Range("E43").Value = '  -4.79%  '
$ws.Range("E44").Value = '  -8.20%  '
$ws.Range("E45").Value = '  -5.74%  '
$ws.Range("D46").Value = "'32.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.20%  '
$ws.Range("D47").Value = "'0.0422"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.42%  '
$ws.Range("E48").Value = '  -3.47%  '
$ws.Range("E49").Value = '  -8.64%  '
$ws.Range("E50").Value = '  -4.31%  '
$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.12%  '

Write-Output "applied 81 cell updates"
